$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.340.57"
$ws.Range("E2").Value = "  +3.11%  "

$ws.Range("D3").Value = "2.060.73"
$ws.Range("E3").Value = "  +5.73%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'235.54"
$ws.Range("E5").Value = "  +3.26%  "

$ws.Range("E6").Value = "  +4.37%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.59"
$ws.Range("E8").Value = "  +9.47%  "

$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  +4.71%  "

$ws.Range("D10").Value = "'57.67"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("D11").Value = "'0.0758"
$ws.Range("E11").Value = "  +3.93%  "

$ws.Range("E12").Value = "  +4.96%  "

$ws.Range("D13").Value = "2.363.25"
$ws.Range("E13").Value = "  +5.70%  "

$ws.Range("D14").Value = "'14.26"
$ws.Range("E14").Value = "  +4.52%  "

$ws.Range("D15").Value = "'20.79"
$ws.Range("E15").Value = "  +7.56%  "

$ws.Range("D16").Value = "'0.772"
$ws.Range("E16").Value = "  +5.05%  "

$ws.Range("D17").Value = "'5.16"
$ws.Range("E17").Value = "  +4.92%  "

$ws.Range("D18").Value = "2.062.14"
$ws.Range("E18").Value = "  +5.76%  "

$ws.Range("D19").Value = "37.238.58"
$ws.Range("E19").Value = "  +3.12%  "

$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "  +24.42%  "

$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").Value = "0.0₃0807"
$ws.Range("E22").Value = "  +3.12%  "

$ws.Range("D23").Value = "'224.39"
$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +6.40%  "

$ws.Range("E26").Value = "  +2.87%  "

$ws.Range("D27").Value = "'163.35"
$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("D28").Value = "'8.81"
$ws.Range("E28").Value = "  +5.75%  "

$ws.Range("E29").Value = "  +10.77%  "

$ws.Range("E30").Value = "  +8.69%  "

$ws.Range("E31").Value = "  +3.19%  "

$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("E33").Value = "  +18.00%  "

$ws.Range("D34").Value = "'0.0625"
$ws.Range("E34").Value = "  +5.11%  "

$ws.Range("E35").Value = "  +4.38%  "

$ws.Range("E36").Value = "  +7.53%  "

$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("E39").Value = "  +7.48%  "

$ws.Range("D40").Value = "'5.83"
$ws.Range("E40").Value = "  +17.52%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'0.0955"
$ws.Range("E42").Value = "  +12.19%  "

$ws.Range("D43").Value = "'4.40"
$ws.Range("E43").Value = "  +29.57%  "

$ws.Range("D44").Value = "1.466.38"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("D45").Value = "'95.07"
$ws.Range("E45").Value = "  +11.37%  "

$ws.Range("E46").Value = "  +6.78%  "

$ws.Range("D47").Value = "'16.07"
$ws.Range("E47").Value = "  +11.05%  "

$ws.Range("E48").Value = "  +5.73%  "

$ws.Range("E49").Value = "  +5.62%  "

$ws.Range("D50").Value = "'7.25"
$ws.Range("E50").Value = "  +9.12%  "

$ws.Range("E51").Value = "  +2.65%  "

